$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 505.21738
$ws.Range("I5").Value = 536
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 536
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -421
$ws.Range("N5").Value = -530
$ws.Range("H15").Value = 797.8570999999999
$ws.Range("I15").Value = 797.8570999999999
$ws.Range("K15").Value = 2393.5713
$ws.Range("M15").Value = -2224.5713
$ws.Range("H38").Value = 5996
$ws.Range("I38").Value = 5776.6665
$ws.Range("J38").Value = 6325
$ws.Range("K38").Value = 17329.9995
$ws.Range("L38").Value = 18975
$ws.Range("M38").Value = -16957.9995
$ws.Range("N38").Value = -19719
$ws.Range("H42").Value = 381.5
$ws.Range("I42").Value = 150
$ws.Range("J42").Value = 497.25
$ws.Range("K42").Value = 450
$ws.Range("L42").Value = 1491.75
$ws.Range("M42").Value = -220
$ws.Range("N42").Value = -1951.75
$ws.Range("H51").Value = 3474.8333
$ws.Range("I51").Value = 2325
$ws.Range("J51").Value = 4049.75
$ws.Range("K51").Value = 2325
$ws.Range("L51").Value = 4049.75
$ws.Range("M51").Value = -1841
$ws.Range("N51").Value = -5017.75
$ws.Range("H76").Value = 10536.579
$ws.Range("I76").Value = 12016.333
$ws.Range("K76").Value = 12016.333
$ws.Range("M76").Value = -11701.333
$ws.Range("H79").Value = 10536.579
$ws.Range("I79").Value = 12016.333
$ws.Range("K79").Value = 12016.333
$ws.Range("M79").Value = -10924.333
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H92").Value = 669.5
$ws.Range("J92").Value = 792.5
$ws.Range("L92").Value = 792.5
$ws.Range("N92").Value = -3288.5
$ws.Range("H132").Value = 2608.9678
$ws.Range("I132").Value = 2137.625
$ws.Range("K132").Value = 6412.875
$ws.Range("M132").Value = -3882.875
$ws.Range("H137").Value = 2555.2083
$ws.Range("I137").Value = 2583.7273
$ws.Range("J137").Value = 2241.5
$ws.Range("K137").Value = 7751.1819
$ws.Range("L137").Value = 6724.5
$ws.Range("M137").Value = -5201.1819
$ws.Range("N137").Value = -11824.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 151.46153
$ws.Range("I5").Value = 121.9
$ws.Range("K5").Value = 121.9
$ws.Range("M5").Value = -9.900000000000006
$ws.Range("H32").Value = 13517912
$ws.Range("I32").Value = 15387293
$ws.Range("J32").Value = 16827.666
$ws.Range("K32").Value = 15387293
$ws.Range("L32").Value = 16827.666
$ws.Range("M32").Value = -15387006
$ws.Range("N32").Value = -17401.666
$ws.Range("H45").Value = 2435.5
$ws.Range("I45").Value = 2087.2
$ws.Range("K45").Value = 2087.2
$ws.Range("M45").Value = -1710.2
$ws.Range("H61").Value = 27781294
$ws.Range("I61").Value = 35716040
$ws.Range("K61").Value = 35716040
$ws.Range("M61").Value = -35715828
$ws.Range("H110").Value = 15918.615
$ws.Range("I110").Value = 17278
$ws.Range("K110").Value = 17278
$ws.Range("M110").Value = -15233
$ws.Range("H113").Value = 80000
$ws.Range("J113").Value = 80000
$ws.Range("L113").Value = 80000
$ws.Range("N113").Value = -88678
$ws.Range("H136").Value = 27781294
$ws.Range("I136").Value = 35716040
$ws.Range("K136").Value = 107148120
$ws.Range("M136").Value = -107145570

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 151.46153
$ws.Range("I4").Value = 121.9
$ws.Range("K4").Value = 121.9
$ws.Range("M4").Value = -6.900000000000006
$ws.Range("H86").Value = 19379.723
$ws.Range("I86").Value = 10411.385
$ws.Range("J86").Value = 42697.4
$ws.Range("K86").Value = 10411.385
$ws.Range("L86").Value = 42697.4
$ws.Range("M86").Value = -9288.385
$ws.Range("N86").Value = -44943.4
$ws.Range("H89").Value = 19379.723
$ws.Range("I89").Value = 10411.385
$ws.Range("J89").Value = 42697.4
$ws.Range("K89").Value = 52056.925
$ws.Range("L89").Value = 213487
$ws.Range("M89").Value = -46440.925
$ws.Range("N89").Value = -224719
$ws.Range("H105").Value = 8469.77
$ws.Range("I105").Value = 9828
$ws.Range("J105").Value = 999.5
$ws.Range("K105").Value = 9828
$ws.Range("L105").Value = 999.5
$ws.Range("M105").Value = -8081
$ws.Range("N105").Value = -4493.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 275.5263
$ws.Range("I7").Value = 117
$ws.Range("J7").Value = 348.69232
$ws.Range("K7").Value = 117
$ws.Range("L7").Value = 348.69232
$ws.Range("M7").Value = -4
$ws.Range("N7").Value = -574.69232
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H58").Value = 2118.5356
$ws.Range("J58").Value = 3634.25
$ws.Range("L58").Value = 3634.25
$ws.Range("N58").Value = -4040.25
$ws.Range("H109").Value = 47883.69
$ws.Range("J109").Value = 46589.91
$ws.Range("L109").Value = 46589.91
$ws.Range("N109").Value = -48669.91
$ws.Range("H136").Value = 2118.5356
$ws.Range("J136").Value = 3634.25
$ws.Range("L136").Value = 10902.75
$ws.Range("N136").Value = -16002.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1434.2609
$ws.Range("I5").Value = 947.6667
$ws.Range("J5").Value = 3186
$ws.Range("K5").Value = 2843.0001
$ws.Range("L5").Value = 9558
$ws.Range("M5").Value = -2731.0001
$ws.Range("N5").Value = -9782
$ws.Range("H107").Value = 552.44446
$ws.Range("J107").Value = 779
$ws.Range("L107").Value = 2337
$ws.Range("N107").Value = -6177
$ws.Range("H113").Value = 2308
$ws.Range("J113").Value = 3238.2
$ws.Range("L113").Value = 9714.599999999999
$ws.Range("N113").Value = -14054.6
$ws.Range("H132").Value = 4446515
$ws.Range("J132").Value = 5557643
$ws.Range("L132").Value = 50018787
$ws.Range("N132").Value = -50023847
$ws.Range("H135").Value = 1434.2609
$ws.Range("I135").Value = 947.6667
$ws.Range("J135").Value = 3186
$ws.Range("K135").Value = 8529.0003
$ws.Range("L135").Value = 28674
$ws.Range("M135").Value = -5994.0003
$ws.Range("N135").Value = -33744
$ws.Range("H140").Value = 1723.8235
$ws.Range("I140").Value = 1172.2307
$ws.Range("K140").Value = 3516.6921
$ws.Range("M140").Value = 1663.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3828.4167
$ws.Range("I80").Value = 2775.3333
$ws.Range("J80").Value = 4881.5
$ws.Range("K80").Value = 2775.3333
$ws.Range("L80").Value = 4881.5
$ws.Range("M80").Value = -1777.3333
$ws.Range("N80").Value = -6877.5
$ws.Range("H83").Value = 3828.4167
$ws.Range("I83").Value = 2775.3333
$ws.Range("J83").Value = 4881.5
$ws.Range("K83").Value = 13876.6665
$ws.Range("L83").Value = 24407.5
$ws.Range("M83").Value = -8884.666499999999
$ws.Range("N83").Value = -34391.5
$ws.Range("H102").Value = 3377.4814
$ws.Range("I102").Value = 2934.6316
$ws.Range("K102").Value = 2934.6316
$ws.Range("M102").Value = -1312.6316

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1333.4584
$ws.Range("I46").Value = 637.19354
$ws.Range("J46").Value = 2603.1177
$ws.Range("K46").Value = 637.19354
$ws.Range("L46").Value = 2603.1177
$ws.Range("M46").Value = -449.19354
$ws.Range("N46").Value = -2979.1177
$ws.Range("H61").Value = 4523.625
$ws.Range("I61").Value = 3787.7
$ws.Range("K61").Value = 3787.7
$ws.Range("M61").Value = -3585.7
$ws.Range("H113").Value = 4523.625
$ws.Range("I113").Value = 3787.7
$ws.Range("K113").Value = 3787.7
$ws.Range("M113").Value = -1617.7
$ws.Range("H122").Value = 4283.61
$ws.Range("I122").Value = 3509.7827
$ws.Range("J122").Value = 5272.3887
$ws.Range("K122").Value = 10529.3481
$ws.Range("L122").Value = 15817.1661
$ws.Range("M122").Value = -8079.348100000001
$ws.Range("N122").Value = -20717.1661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 589.625
$ws.Range("I81").Value = 620.5
$ws.Range("J81").Value = 497
$ws.Range("K81").Value = 1241
$ws.Range("L81").Value = 994
$ws.Range("M81").Value = -180
$ws.Range("N81").Value = -3116
$ws.Range("H84").Value = 589.625
$ws.Range("I84").Value = 620.5
$ws.Range("J84").Value = 497
$ws.Range("K84").Value = 6205
$ws.Range("L84").Value = 4970
$ws.Range("M84").Value = -901
$ws.Range("N84").Value = -15578
$ws.Range("H132").Value = 5344.115
$ws.Range("I132").Value = 5357.88
$ws.Range("K132").Value = 16073.64
$ws.Range("M132").Value = -13543.64
$ws.Range("H136").Value = 1937.08
$ws.Range("I136").Value = 1716.9474
$ws.Range("K136").Value = 5150.8422
$ws.Range("M136").Value = -2600.8422
